$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: shared-string header, index shifts from 10 -> 11 (new si inserted) but text unchanged
$ws.Range("A1").Value = "HK_R_acc_G"

# A2:A49 recomputed result values
$ws.Range("A2").Value = 76.2471395881007
$ws.Range("A3").Value = 75.92677345537757
$ws.Range("A4").Value = 74.0045766590389
$ws.Range("A5").Value = 74.59954233409611
$ws.Range("A6").Value = 74.64530892448512
$ws.Range("A7").Value = 75.92677345537757
$ws.Range("A8").Value = 74.55377574370709
$ws.Range("A9").Value = 75.51487414187643
$ws.Range("A10").Value = 74.41647597254004
$ws.Range("A11").Value = 74.41647597254004
$ws.Range("A12").Value = 75.28604118993135
$ws.Range("A13").Value = 75.74370709382151
$ws.Range("A14").Value = 74.55377574370709
$ws.Range("A15").Value = 75.05720823798627
$ws.Range("A16").Value = 74.8283752860412
$ws.Range("A17").Value = 74.50800915331808
$ws.Range("A18").Value = 75.78947368421053
$ws.Range("A19").Value = 75.83524027459954
$ws.Range("A20").Value = 77.57437070938215
$ws.Range("A21").Value = 75.33180778032037
$ws.Range("A22").Value = 77.7116704805492
$ws.Range("A23").Value = 74.96567505720823
$ws.Range("A24").Value = 74.69107551487414
$ws.Range("A25").Value = 74.8283752860412
$ws.Range("A26").Value = 74.59954233409611
$ws.Range("A27").Value = 74.55377574370709
$ws.Range("A28").Value = 74.279176201373
$ws.Range("A29").Value = 75.42334096109839
$ws.Range("A30").Value = 74.59954233409611
$ws.Range("A31").Value = 74.87414187643022
$ws.Range("A32").Value = 74.0045766590389
$ws.Range("A33").Value = 74.23340961098398
$ws.Range("A34").Value = 74.0045766590389
$ws.Range("A35").Value = 74.87414187643022
$ws.Range("A36").Value = 75.1029748283753
$ws.Range("A37").Value = 81.09839816933638
$ws.Range("A38").Value = 74.41647597254004
$ws.Range("A39").Value = 74.8283752860412
$ws.Range("A40").Value = 75.74370709382151
$ws.Range("A41").Value = 74.78260869565217
$ws.Range("A42").Value = 75.24027459954233
$ws.Range("A43").Value = 74.96567505720823
$ws.Range("A44").Value = 75.01144164759725
$ws.Range("A45").Value = 75.33180778032037
$ws.Range("A46").Value = 74.0045766590389
$ws.Range("A47").Value = 74.279176201373
$ws.Range("A48").Value = 76.10983981693363
$ws.Range("A49").Value = 74.64530892448512
